$d = $word.ActiveDocument

# Mapping of old text to new text, applied in document order.
# (Doc order matters: "24÷8=" is both an existing cell value and the
# result of another cell's replacement, so we must replace the original
# "24÷8=" cell before creating a new "24÷8=" elsewhere.)
$pairs = @(
    @("2024-06-24 Monday", "2024-06-25 Tuesday"),
    @("28÷5=", "38÷3="),
    @("77÷9=", "29÷9="),
    @("62÷9=", "56÷6="),
    @("92÷2=", "66÷7="),
    @("20÷9=", "67÷5="),
    @("34÷7=", "48÷5="),
    @("41÷5=", "48÷4="),
    @("53÷3=", "67÷2="),
    @("99÷7=", "99÷9="),
    @("51÷8=", "90÷7="),
    @("59÷7=", "40÷9="),
    @("39÷2=", "63÷9="),
    @("45÷4=", "93÷3="),
    @("29÷5=", "96÷5="),
    @("24÷8=", "55÷5="),
    @("67÷3=", "19÷2="),
    @("93÷5=", "26÷6="),
    @("71÷5=", "24÷8="),
    @("26÷7=", "44÷5="),
    @("41÷9=", "49÷7="),
    @("85÷8=", "86÷6="),
    @("53÷9=", "96÷6="),
    @("69÷7=", "56÷3="),
    @("36÷5=", "82÷7="),
    @("48÷2=", "17÷6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
